# Apply the "Added code for H2H resolution" edit to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Calculation mode: workbook was "manual", make it automatic again.
$excel.Calculation = -4105   # xlCalculationAutomatic

# 2. New column Q: header "XTM" in Q1, formatted like the other headers
#    (font/fill copied from P1) but with a thin left/right-only border.
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q1").Value = "XTM"
$ws.Range("Q1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft  -> thin
$ws.Range("Q1").Borders.Item(10).LineStyle = 1  # xlEdgeRight -> thin
$ws.Range("Q1").Borders.Item(8).LineStyle = -4142  # xlEdgeTop    -> none
$ws.Range("Q1").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> none

# 3. New data cell Q11 (row 11 is a "Girls" row), formatted like P11
#    but with the same thin left/right-only border as Q1.
$ws.Range("P11").Copy()
$ws.Range("Q11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q11").Value = 1
$ws.Range("Q11").Borders.Item(7).LineStyle = 1
$ws.Range("Q11").Borders.Item(10).LineStyle = 1
$ws.Range("Q11").Borders.Item(8).LineStyle = -4142
$ws.Range("Q11").Borders.Item(9).LineStyle = -4142

# 4. Score corrections (H2H resolution).
$ws.Range("E6").Value = 2
$ws.Range("O11").Value = 1
$ws.Range("N13").Value = 1

# 5. Update the active selection to match the author's final cursor position.
$ws.Range("O14").Select()

$wb.Save()
